$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 362; existing rows 362-374 shift down to 364-376.
$ws.Rows("362:363").Insert()

# New row 362 data
$ws.Cells.Item(362, 1).Value = 5
$ws.Cells.Item(362, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(362, 3).Value = "Maule"
$ws.Cells.Item(362, 4).Value = 44509
$ws.Cells.Item(362, 5).Value = 7
$ws.Cells.Item(362, 6).Value = 100112002
$ws.Cells.Item(362, 7).Value = "Pimiento"
$ws.Cells.Item(362, 8).Value = "Cuatro cascos verde"
$ws.Cells.Item(362, 9).Value = "Primera"
$ws.Cells.Item(362, 10).Value = 200
$ws.Cells.Item(362, 11).Value = 27000
$ws.Cells.Item(362, 12).Value = 27000
$ws.Cells.Item(362, 13).Value = 27000
$ws.Cells.Item(362, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(362, 15).Value = "Región del Maule"
$ws.Cells.Item(362, 16).Value = 1800
$ws.Cells.Item(362, 17).Value = 15
$ws.Cells.Item(362, 18).Value = "Hortaliza"

# New row 363 data
$ws.Cells.Item(363, 1).Value = 5
$ws.Cells.Item(363, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(363, 3).Value = "Maule"
$ws.Cells.Item(363, 4).Value = 44509
$ws.Cells.Item(363, 5).Value = 7
$ws.Cells.Item(363, 6).Value = 100112002
$ws.Cells.Item(363, 7).Value = "Pimiento"
$ws.Cells.Item(363, 8).Value = "Zafiro rojo"
$ws.Cells.Item(363, 9).Value = "Primera"
$ws.Cells.Item(363, 10).Value = 200
$ws.Cells.Item(363, 11).Value = 45000
$ws.Cells.Item(363, 12).Value = 45000
$ws.Cells.Item(363, 13).Value = 45000
$ws.Cells.Item(363, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(363, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(363, 16).Value = 3000
$ws.Cells.Item(363, 17).Value = 15
$ws.Cells.Item(363, 18).Value = "Hortaliza"

"done"
